$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value pairs derived from the source diff.
# Values that look like plain numbers are prefixed with a leading apostrophe
# so Excel stores them as text (matching the workbook's original inlineStr
# typing) instead of silently converting them to numeric values and
# dropping formatting such as trailing zeros.
$updates = @(
    @("D2", "69.619.13"),
    @("E2", "  -1.48%  "),
    @("D3", "3.528.81"),
    @("E3", "  -1.66%  "),
    @("D4", "'1.00"),
    @("E4", "  +0.03%  "),
    @("D5", "'610.03"),
    @("E5", "  +3.51%  "),
    @("D6", "'183.59"),
    @("E6", "  -1.94%  "),
    @("D7", "'0.614"),
    @("E8", "  -0.05%  "),
    @("D9", "'0.214"),
    @("E9", "  +5.12%  "),
    @("D10", "'0.639"),
    @("E10", "  -1.81%  "),
    @("D11", "'53.21"),
    @("E11", "  -2.42%  "),
    @("D12", "'0.0000306"),
    @("E12", "  -1.67%  "),
    @("D13", "'9.41"),
    @("E13", "  -1.50%  "),
    @("D14", "4.091.36"),
    @("E14", "  -1.49%  "),
    @("D15", "'593.42"),
    @("E15", "  +5.61%  "),
    @("B16", "WrappedEther"),
    @("C16", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D16", "3.579.71"),
    @("E16", "  +1.31%  "),
    @("B17", "WrappedBTC"),
    @("C17", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"),
    @("D17", "69.710.66"),
    @("E17", "  -1.27%  "),
    @("D18", "'12.59"),
    @("E18", "  +0.86%  "),
    @("D19", "'18.80"),
    @("E19", "  -4.76%  "),
    @("E20", "  -0.49%  "),
    @("D21", "'0.986"),
    @("E21", "  -3.18%  "),
    @("D22", "'17.41"),
    @("E22", "  -3.15%  "),
    @("D23", "'4.71"),
    @("E23", "  +0.66%  "),
    @("D24", "'98.68"),
    @("E24", "  +3.22%  "),
    @("E25", "  -1.31%  "),
    @("D26", "'2.95"),
    @("E26", "  -0.93%  "),
    @("E27", "  -6.07%  "),
    @("D28", "'9.60"),
    @("E28", "  +4.46%  "),
    @("D29", "'32.11"),
    @("E29", "  -0.53%  "),
    @("D30", "'6.97"),
    @("E30", "  -6.00%  "),
    @("E31", "  -3.65%  "),
    @("E32", "  -1.80%  "),
    @("D33", "'63.32"),
    @("E33", "  -2.89%  "),
    @("B34", "Fetch.AI"),
    @("C34", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"),
    @("D34", "'3.22"),
    @("E34", "  -4.26%  "),
    @("B35", "dogwifhat"),
    @("C35", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"),
    @("D35", "'3.54"),
    @("E35", "  +15.80%  "),
    @("D36", "'529.67"),
    @("E36", "  -5.84%  "),
    @("D38", "'0.399"),
    @("E38", "  -5.51%  "),
    @("D39", "'36.89"),
    @("E39", "  -2.97%  "),
    @("D40", "3.530.08"),
    @("E40", "  +5.46%  "),
    @("D41", "0.0₃0774"),
    @("E41", "  -0.49%  "),
    @("E42", "  +3.94%  "),
    @("E43", "  +0.27%  "),
    @("D44", "'0.0454"),
    @("E44", "  +1.69%  "),
    @("D45", "'2.91"),
    @("E45", "  -2.37%  "),
    @("D46", "'3.40"),
    @("E46", "  -3.78%  "),
    @("E47", "  +2.83%  "),
    @("D48", "'9.07"),
    @("E48", "  -4.06%  "),
    @("E49", "  +0.40%  "),
    @("E50", "  -6.43%  "),
    @("D51", "'135.07"),
    @("E51", "  -2.11%  ")
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $ws.Range($cellRef).Value = $newValue
}
